$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Column L (Seasonality Index) updates
$ws.Range("L2").Value = 0.91
$ws.Range("L3").Value = 0.83
$ws.Range("L4").Value = 1.17
$ws.Range("L5").Value = 1.05
$ws.Range("L6").Value = 1.1
$ws.Range("L7").Value = 1.07
$ws.Range("L8").Value = 1.04
$ws.Range("L9").Value = 1.03
$ws.Range("L10").Value = 0.9
$ws.Range("L11").Value = 0.98
$ws.Range("L12").Value = 1.03
$ws.Range("L13").Value = 0.91
$ws.Range("L14").Value = 0.88
$ws.Range("L15").Value = 0.84
$ws.Range("L16").Value = 1.11
$ws.Range("L17").Value = 0.85

# Column H (Inventory Coverage) updates
$ws.Range("H4").Value = 355.14
$ws.Range("H5").Value = 354.14
$ws.Range("H6").Value = 353.14
$ws.Range("H7").Value = 352.14
$ws.Range("H8").Value = 351.14
$ws.Range("H9").Value = 350.14
$ws.Range("H10").Value = 349.14
$ws.Range("H11").Value = 348.14
$ws.Range("H12").Value = 347.14
$ws.Range("H13").Value = 346.14
$ws.Range("H14").Value = 345.14
$ws.Range("H15").Value = 344.14
